$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking Price values so Excel
# does not silently convert them to floating point numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values from the source diff.
$ws.Range("D2").Value2 = "65.784.30"
$ws.Range("E2").Value2 = "  +0.23%  "
$ws.Range("D3").Value2 = "2.680.51"
$ws.Range("E4").Value2 = "  -0.02%  "
$ws.Range("D5").Value2 = "601.83"
$ws.Range("E5").Value2 = "  -1.17%  "
$ws.Range("D6").Value2 = "157.05"
$ws.Range("E6").Value2 = "  -0.34%  "
$ws.Range("E7").Value2 = "  +0.01%  "
$ws.Range("D8").Value2 = "0.625"
$ws.Range("E8").Value2 = "  +6.43%  "
$ws.Range("D9").Value2 = "0.131"
$ws.Range("E9").Value2 = "  +5.34%  "
$ws.Range("E10").Value2 = "  -0.06%  "
$ws.Range("E11").Value2 = "  -3.27%  "
$ws.Range("E12").Value2 = "  -0.17%  "
$ws.Range("D13").Value2 = "29.44"
$ws.Range("E13").Value2 = "  -2.61%  "
$ws.Range("E14").Value2 = "  -1.45%  "
$ws.Range("D15").Value2 = "3.162.68"
$ws.Range("E15").Value2 = "  -0.36%  "
$ws.Range("D16").Value2 = "65.635.20"
$ws.Range("E16").Value2 = "  +0.24%  "
$ws.Range("D17").Value2 = "2.683.28"
$ws.Range("E17").Value2 = "  -0.51%  "
$ws.Range("E18").Value2 = "  +1.47%  "
$ws.Range("D19").Value2 = "4.82"
$ws.Range("E19").Value2 = "  -1.32%  "
$ws.Range("D20").Value2 = "7.60"
$ws.Range("E20").Value2 = "  +0.92%  "
$ws.Range("D21").Value2 = "352.62"
$ws.Range("E21").Value2 = "  -1.92%  "
$ws.Range("D23").Value2 = "69.75"
$ws.Range("E23").Value2 = "  -0.53%  "
$ws.Range("E24").Value2 = "  +4.78%  "
$ws.Range("D25").Value2 = "9.68"
$ws.Range("E25").Value2 = "  -1.38%  "
$ws.Range("E26").Value2 = "  +1.06%  "
$ws.Range("E27").Value2 = "  -1.01%  "
$ws.Range("E28").Value2 = "  -5.44%  "
$ws.Range("D29").Value2 = "8.10"
$ws.Range("E29").Value2 = "  -1.96%  "
$ws.Range("E30").Value2 = "  +0.37%  "
$ws.Range("B31").Value2 = "PancakeSwap"
$ws.Range("C31").Value2 = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value2 = "2.15"
$ws.Range("E31").Value2 = "  -2.55%  "
$ws.Range("B32").Value2 = "Bittensor"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value2 = "531.81"
$ws.Range("E32").Value2 = "  -0.50%  "
$ws.Range("E33").Value2 = "  -1.51%  "
$ws.Range("D34").Value2 = "6.49"
$ws.Range("E34").Value2 = "  -2.71%  "
$ws.Range("E35").Value2 = "  +0.92%  "
$ws.Range("E36").Value2 = "  -1.22%  "
$ws.Range("E37").Value2 = "  -0.77%  "
$ws.Range("E38").Value2 = "  +0.01%  "
$ws.Range("D39").Value2 = "158.29"
$ws.Range("E39").Value2 = "  -2.79%  "
$ws.Range("E40").Value2 = "  -2.36%  "
$ws.Range("E41").Value2 = "  +0.02%  "
$ws.Range("D42").Value2 = "164.94"
$ws.Range("E42").Value2 = "  -2.22%  "
$ws.Range("E43").Value2 = "  -0.63%  "
$ws.Range("D44").Value2 = "2.33"
$ws.Range("E44").Value2 = "  +2.86%  "
$ws.Range("E45").Value2 = "  -0.47%  "
$ws.Range("D46").Value2 = "22.95"
$ws.Range("E46").Value2 = "  -2.25%  "
$ws.Range("E47").Value2 = "  -2.03%  "
$ws.Range("E48").Value2 = "  -2.70%  "
$ws.Range("E49").Value2 = "  +13.68%  "
$ws.Range("E50").Value2 = "  +2.65%  "
$ws.Range("D51").Value2 = "20.17"
$ws.Range("E51").Value2 = "  -4.22%  "
